$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the old rows 4-10 entirely - only the two newest listings remain.
$ws.Range("A4:H10").EntireRow.Delete()

# Row 2 becomes the "エクセルでの在庫管理システム構築依頼" listing (previously row 8).
$ws.Range("A2").Value2 = "2025-12-30 06:30:09"
$ws.Range("B2").Value2 = "【急募】エクセルでの在庫管理システム構築依頼"
$ws.Range("C2").Value2 = "システム開発"
$ws.Range("D2").Value2 = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E2").Value2 = "期限情報なし"
$ws.Range("F2").Value2 = "https://www.lancers.jp/work/detail/5463183"
$ws.Range("G2").Value2 = 45
$ws.Range("H2").Value2 = "◇管理"

# Row 3 becomes the new "テレグラム風メッセージアプリ" listing.
$ws.Range("A3").Value2 = "2025-12-30 06:30:09"
$ws.Range("B3").Value2 = "【急募】テレグラム風メッセージアプリのバグ修正と機能追加"
$ws.Range("C3").Value2 = "システム開発"
$ws.Range("D3").Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value2 = "期限情報なし"
$ws.Range("F3").Value2 = "https://www.lancers.jp/work/detail/5463296"
$ws.Range("G3").Value2 = 38
$ws.Range("H3").Value2 = "◇アプリ"

# The worksheet-level Hyperlinks collection is rebuilt wholesale here since the
# old rows (and their hyperlinks) are gone; re-create the two that remain.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5463183") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5463296") | Out-Null

# Make sure F2/F3 keep the underlined "Hyperlink" look after re-adding the links.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

# Narrow columns B and D to their new widths (21/51/12/39/... -> 21/30/12/28/...).
$ws.Columns.Item(2).ColumnWidth = 29.17
$ws.Columns.Item(4).ColumnWidth = 27.17
